{"js": "// Update benchmark stats table for the DaCapo / Shenandoah GC (avrora, heap-4G) doc.\n// The document body contains a single table, one value per row (column 1).\n// Rows 0-2 become the rollup \"0M\" placeholder, row 3 is corrected to 35,\n// rows 5-11 (the per-run detail stats) are refreshed with new measurements,\n// and the final three \"raw run\" rows (43-45), which previously held the full\n// tab-separated dump of a run, are collapsed back down to just their\n// single summary value (matching rows 0-2's original values).\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nconst updates = [\n  [0, \"0M\"],\n  [1, \"0M\"],\n  [2, \"0M\"],\n  [3, \"35\"],\n  [5, \"0.00101\"],\n  [6, \"0.00022\"],\n  [7, \"0.00007\"],\n  [8, \"0.00028\"],\n  [9, \"0.00043\"],\n  [10, \"0.00059\"],\n  [11, \"0.00936\"],\n  [43, \"99.99\"],\n  [44, \"0.01\"],\n  [45, \"65\"],\n];\n\nfor (const [rowIndex, newText] of updates) {\n  const cell = table.getCell(rowIndex, 0);\n  cell.value = newText;\n}\n\nawait context.sync();\n", "ps1": "# Update benchmark stats table for the DaCapo / Shenandoah GC (avrora, heap-4G) doc.\n# The document is a single table, one value per row (column 1).\n# Rows 1-3 become the rollup \"0M\" placeholder, row 4 is corrected to 35,\n# rows 6-12 (the per-run detail stats) are refreshed with new measurements,\n# and the final three \"raw run\" rows (44-46), which previously held the full\n# tab-separated dump of a run, are collapsed back down to just their\n# single summary value (matching rows 1-3's original values).\n\n$d = $word.ActiveDocument\n$t = $d.Tables(1)\n\nfunction Set-CellText($table, $row, $text) {\n    $cell = $table.Cell($row, 1)\n    $cell.Range.Text = $text\n}\n\nSet-CellText $t 1 \"0M\"\nSet-CellText $t 2 \"0M\"\nSet-CellText $t 3 \"0M\"\nSet-CellText $t 4 \"35\"\n\nSet-CellText $t 6 \"0.00101\"\nSet-CellText $t 7 \"0.00022\"\nSet-CellText $t 8 \"0.00007\"\nSet-CellText $t 9 \"0.00028\"\nSet-CellText $t 10 \"0.00043\"\nSet-CellText $t 11 \"0.00059\"\nSet-CellText $t 12 \"0.00936\"\n\nSet-CellText $t 44 \"99.99\"\nSet-CellText $t 45 \"0.01\"\nSet-CellText $t 46 \"65\"\n"}
